$d = $word.ActiveDocument

# Step 1: remove the existing _GoBack bookmark (it will be re-added at the new location below)
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

# Step 2: replace paragraph 10 ("All of the number columns...") together with the blank
# separator paragraph 11 that follows it, with four new paragraphs:
#   - new explanatory text about Unemployment/Crime towns (ending with the _GoBack bookmark)
#   - an empty separator paragraph
#   - the original "All of the number columns..." paragraph, relocated
#   - another empty separator paragraph
$p10 = $d.Paragraphs.Item(10)
$p11 = $d.Paragraphs.Item(11)
$rng = $d.Range($p10.Range.Start, $p11.Range.End)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">The Unemployment data had additional towns that the Crime data did not.  The towns that did not exist in the crime </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>csv</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> were removed from the merged data using Google fusion tables.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>All of the number columns except population were changed to have one decimal place to keep the data consistent.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">  This change was made </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve">with a combination of </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>sed</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t xml:space="preserve"> and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>awk</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/>
          <w:sz w:val="22"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    </w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng.InsertXML($xml)

Write-Host "Paragraph count after edit:" $d.Paragraphs.Count
